# Plantilla_Resultados_ESAComp - fill in the results obtained from ESAComp
# ("Valor extraído de ESACOMP") on Hoja1, column C, and let the existing
# error-percentage formulas underneath (dE1/dE2/dS1/dS2, rows 12-15) use
# them. C14/C15 did not have their formulas yet, so add them too (mirrors
# the ones already present in C12/C13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# APARTADO 1
$ws.Range("C3").Value = 276      # Ef  [GPa]
$ws.Range("C4").Value = 5654     # Sf  [MPa]
$ws.Range("C5").Value = 4.67     # Em  [GPa]
$ws.Range("C6").Value = 121      # Sm  [MPa]

$ws.Range("C8").Value = 161.22             # E1 [GPa]
$ws.Range("C9").Value = 10.856             # E2 [GPa]
$ws.Range("C10").Value = 3313.5410000000002 # S1 [MPa]
$ws.Range("C11").Value = 277.93799999999999 # S2 [MPa]

# dS1 / dS2 formulas (same pattern as the existing dE1/dE2 ones in C12/C13)
$ws.Range("C14").Formula = "=(C10-2724)*100/2724"
$ws.Range("C15").Formula = "=(C11-111)*100/111"

$excel.Calculate()

# Leave the selection where the author ended up after entering the data
$ws.Range("C24").Select()
